# Refresh the cryptos list: updated prices / %-change figures, plus three
# pairs of adjacent rows whose coin (name, link, price, volume) got swapped
# (InternetComputer<->RenderToken, Hedera<->Cosmos, PEPE<->TheGraph).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text formatting on every cell we touch first, so numeric-looking
# strings (e.g. "0.719", "0.0000291") are written as text, not coerced into
# Excel numbers (which would mangle values like "0.0000291" -> 2.91E-05).
$cells = @(
    "D2", "E2", "D3", "E3", "E4", "D5", "E5", "D6",
    "E6", "D7", "E7", "E8", "D9", "E9", "E10", "D11",
    "E11", "D12", "E12", "D13", "E13", "D14", "E14", "D15",
    "E15", "D16", "E16", "D17", "E17", "D18", "E18", "E19",
    "D20", "E20", "D21", "E21", "D22", "E22", "D23", "E23",
    "E24", "B25", "C25", "D25", "E25", "B26", "C26", "D26",
    "E26", "E27", "E28", "D29", "E29", "D30", "E30", "E31",
    "B32", "C32", "D32", "E32", "B33", "C33", "D33", "E33",
    "D34", "E34", "D35", "E35", "D36", "E36", "B37", "C37",
    "D37", "E37", "B38", "C38", "D38", "E38", "E39", "E40",
    "D41", "E41", "E42", "D43", "E43", "D44", "E44", "E45",
    "D46", "E46", "D47", "E47", "D48", "E48", "E49", "E50",
    "D51", "E51"
)
foreach ($addr in $cells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "69.035.13"
$ws.Range("E2").Value = "  +2.05%  "
$ws.Range("D3").Value = "3.728.64"
$ws.Range("E3").Value = "  +1.44%  "
$ws.Range("E4").Value = "  -0.52%  "
$ws.Range("D5").Value = "617.20"
$ws.Range("E5").Value = "  +6.61%  "
$ws.Range("D6").Value = "188.14"
$ws.Range("E6").Value = "  +7.85%  "
$ws.Range("D7").Value = "0.640"
$ws.Range("E7").Value = "  +1.68%  "
$ws.Range("E8").Value = "  +0.17%  "
$ws.Range("D9").Value = "0.719"
$ws.Range("E9").Value = "  +1.10%  "
$ws.Range("E10").Value = "  -1.89%  "
$ws.Range("D11").Value = "57.21"
$ws.Range("E11").Value = "  +9.39%  "
$ws.Range("D12").Value = "0.0000291"
$ws.Range("E12").Value = "  -2.33%  "
$ws.Range("D13").Value = "10.62"
$ws.Range("E13").Value = "  +0.01%  "
$ws.Range("D14").Value = "4.310.01"
$ws.Range("E14").Value = "  +1.18%  "
$ws.Range("D15").Value = "3.720.12"
$ws.Range("E15").Value = "  +1.01%  "
$ws.Range("D16").Value = "19.40"
$ws.Range("E16").Value = "  +0.38%  "
$ws.Range("D17").Value = "13.03"
$ws.Range("E17").Value = "  +0.51%  "
$ws.Range("D18").Value = "0.127"
$ws.Range("E18").Value = "  -0.50%  "
$ws.Range("E19").Value = "  +0.35%  "
$ws.Range("D20").Value = "68.832.68"
$ws.Range("E20").Value = "  +1.80%  "
$ws.Range("D21").Value = "412.10"
$ws.Range("E21").Value = "  +1.18%  "
$ws.Range("D22").Value = "4.62"
$ws.Range("E22").Value = "  +2.04%  "
$ws.Range("D23").Value = "89.48"
$ws.Range("E23").Value = "  +1.68%  "
$ws.Range("E24").Value = "  -0.29%  "
$ws.Range("B25").Value = "RenderToken"
$ws.Range("C25").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D25").Value = "11.19"
$ws.Range("E25").Value = "  +4.33%  "
$ws.Range("B26").Value = "InternetComputer(DFINITY)"
$ws.Range("C26").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D26").Value = "12.92"
$ws.Range("E26").Value = "  +1.25%  "
$ws.Range("E27").Value = "  +1.33%  "
$ws.Range("E28").Value = "  -1.32%  "
$ws.Range("D29").Value = "9.67"
$ws.Range("E29").Value = "  +1.80%  "
$ws.Range("D30").Value = "33.29"
$ws.Range("E30").Value = "  +1.50%  "
$ws.Range("E31").Value = "  -10.10%  "
$ws.Range("B32").Value = "Cosmos"
$ws.Range("C32").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D32").Value = "12.66"
$ws.Range("E32").Value = "  -0.13%  "
$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").Value = "0.123"
$ws.Range("E33").Value = "  +4.78%  "
$ws.Range("D34").Value = "627.58"
$ws.Range("E34").Value = "  +6.88%  "
$ws.Range("D35").Value = "44.74"
$ws.Range("E35").Value = "  +0.57%  "
$ws.Range("D36").Value = "66.12"
$ws.Range("E36").Value = "  +0.72%  "
$ws.Range("B37").Value = "TheGraph"
$ws.Range("C37").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D37").Value = "0.418"
$ws.Range("E37").Value = "  +3.87%  "
$ws.Range("B38").Value = "PEPE"
$ws.Range("C38").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D38").Value = "0.0₃0835"
$ws.Range("E38").Value = "  -10.04%  "
$ws.Range("E39").Value = "  -0.14%  "
$ws.Range("E40").Value = "  -0.56%  "
$ws.Range("D41").Value = "0.141"
$ws.Range("E41").Value = "  +3.61%  "
$ws.Range("E42").Value = "  +0.29%  "
$ws.Range("D43").Value = "0.0446"
$ws.Range("E43").Value = "  +1.27%  "
$ws.Range("D44").Value = "2.63"
$ws.Range("E44").Value = "  +1.79%  "
$ws.Range("E45").Value = "  +4.24%  "
$ws.Range("D46").Value = "2.856.43"
$ws.Range("E46").Value = "  +4.57%  "
$ws.Range("D47").Value = "2.76"
$ws.Range("E47").Value = "  +3.68%  "
$ws.Range("D48").Value = "9.13"
$ws.Range("E48").Value = "  -3.47%  "
$ws.Range("E49").Value = "  -1.20%  "
$ws.Range("E50").Value = "  -18.09%  "
$ws.Range("D51").Value = "142.03"
$ws.Range("E51").Value = "  +1.08%  "
